# Lasso regression + normalization + lag1 accuracy results updated:
# Insert a new row for "Lasso Regression+normalization+ lag1 +PCA(2)" right
# after the existing "...+ lag1" row (row 6), pushing the remaining model
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7 (shifts rows 7..11 down to 8..12).
$ws.Rows.Item(7).Insert()

# Copy formatting from the row above (row 6) so the new row matches the
# border/style used by the other interior data rows.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)

# Populate the new row with the PCA(2) result.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Lasso Regression+normalization+ lag1 +PCA(2)"
$ws.Range("C7").Value = 85.044508627085506

# Widen column B to fit the longer model names, and update the active
# selection to match the author's final cursor position.
$ws.Columns.Item(2).ColumnWidth = 45.8
[void]$ws.Range("E4").Select()
